# Bold the "M. Grossi" author name within the 2015 SMAST technical-report
# citation (the "Pre- and Post-Mission Glider CTD Comparison Measurements:
# 19 June 2014 and 6 February 2015 ... SMAST-15-06-01" reference), matching
# the formatting already used on the sibling 2016 citation just above it.

$d = $word.ActiveDocument

$target = "M. Grossi"

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text

    if ($text -like "*SMAST-15-06-01*" -and $text -like "*W.S. Brown and M. Grossi*") {
        $paraStart = $p.Range.Start
        $idx = $text.IndexOf($target)

        if ($idx -ge 0) {
            $boldStart = $paraStart + $idx
            $boldEnd = $boldStart + $target.Length

            # Apply bold only to the non-empty "M. Grossi" sub-range so the
            # surrounding text keeps its existing (non-bold) formatting.
            $boldRange = $d.Range($boldStart, $boldEnd)
            $boldRange.Bold = 1
            $boldRange.Font.Bold = 1
        }
    }
}
